# Week 13 logging update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly logged per-play yardage figures to the
# running game logs stored as space separated strings.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Text + " 22 1 3 1 16 2 6 2 2 4 1 9 5 -4 -3 13 4 1 0 -6 7 2 3 8 1 6 0 4 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Text + " 23 -1 13 6 7 4 7 6 5 6 5 6 0 10 35 10 3 4 10 9 10 12 6"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Text + " 3 7 1 3 0 17 7 3 15 -1 6 2 1 1 10 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Text + " 1 8 9 3 10 5 3 9 7 5 0 18 10 -5 5 2 15 9 34 7 3 28 14 13 16 2 11 7"

# ---------------------------------------------------------------------
# OFF sheet: updated Home/Road offensive season totals.
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 150
$offWs.Range("F2").Value = 58
$offWs.Range("G2").Value = 47
$offWs.Range("H2").Value = 4
$offWs.Range("I2").Value = 6
$offWs.Range("J2").Value = 33
$offWs.Range("N2").Value = 15

$offWs.Range("B3").Value = 10
$offWs.Range("C3").Value = 142
$offWs.Range("E3").Value = 25
$offWs.Range("F3").Value = 63
$offWs.Range("H3").Value = 18
$offWs.Range("I3").Value = 38
$offWs.Range("J3").Value = 40
$offWs.Range("L3").Value = 184
$offWs.Range("M3").Value = 125
$offWs.Range("Q3").Value = 376

# ---------------------------------------------------------------------
# DEF sheet: updated Home/Road defensive season totals.
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 119
$defWs.Range("D2").Value = 6
$defWs.Range("E2").Value = 6
$defWs.Range("F2").Value = 47
$defWs.Range("G2").Value = 22
$defWs.Range("N2").Value = 12

$defWs.Range("C3").Value = 148
$defWs.Range("E3").Value = 21
$defWs.Range("F3").Value = 82
$defWs.Range("G3").Value = 27
$defWs.Range("H3").Value = 16
$defWs.Range("I3").Value = 39
$defWs.Range("J3").Value = 53
$defWs.Range("L3").Value = 212
$defWs.Range("M3").Value = 152
$defWs.Range("Q3").Value = 352

# ---------------------------------------------------------------------
# ST sheet: updated special teams totals plus appended kickoff / return
# game logs.
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 57
$stWs.Range("D2").Value = 36
$stWs.Range("F2").Value = 93
$stWs.Range("G2").Value = 82
$stWs.Range("L2").Value = 27
$stWs.Range("M2").Value = 21

$stWs.Range("B3").Value = 39

$stWs.Range("B4").Value = $stWs.Range("B4").Text + " 60 64"
$stWs.Range("D4").Value = $stWs.Range("D4").Text + " 6 0 0 0"

$stWs.Range("B5").Value = $stWs.Range("B5").Text + " 19 24"
$stWs.Range("D5").Value = $stWs.Range("D5").Text + " 0 16 0 13"

$stWs.Range("D3").Value = $stWs.Range("D3").Text + " 54 47 43 53"

# ---------------------------------------------------------------------
# TURNS sheet: updated Road turnover totals.
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 6
$turnsWs.Range("E3").Value = 5

# ---------------------------------------------------------------------
# PEN sheet: updated penalty totals.
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 13
$penWs.Range("B3").Value = 12
$penWs.Range("D3").Value = 3
$penWs.Range("D4").Value = 8
